$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 96-97, shifting the existing rows 96-110 down to 98-112.
$ws.Rows("96:97").Insert()

# --- New row 96 ---
$ws.Range("A96").Value = 1
$ws.Range("B96").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C96").Value = "Arica y Parinacota"
$ws.Range("D96").Value = 44769
$ws.Range("E96").Value = 15
$ws.Range("F96").Value = "Fruta"
$ws.Range("G96").Value = 100102
$ws.Range("H96").Value = "Cítricos"
$ws.Range("I96").Value = 100102004
$ws.Range("J96").Value = "Mandarina"
$ws.Range("K96").Value = "Clemenuless"
$ws.Range("L96").Value = "Segunda"
$ws.Range("M96").Value = 300
$ws.Range("N96").Value = 17000
$ws.Range("O96").Value = 18000
$ws.Range("P96").Value = 17500
$ws.Range("Q96").Value = "$/caja 20 kilos"
$ws.Range("R96").Value = "Región de Coquimbo"
$ws.Range("S96").Value = 875
$ws.Range("T96").Value = 20

# --- New row 97 ---
$ws.Range("A97").Value = 1
$ws.Range("B97").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C97").Value = "Arica y Parinacota"
$ws.Range("D97").Value = 44769
$ws.Range("E97").Value = 15
$ws.Range("F97").Value = "Fruta"
$ws.Range("G97").Value = 100102
$ws.Range("H97").Value = "Cítricos"
$ws.Range("I97").Value = 100102004
$ws.Range("J97").Value = "Mandarina"
$ws.Range("K97").Value = "Murcott"
$ws.Range("L97").Value = "Tercera"
$ws.Range("M97").Value = 250
$ws.Range("N97").Value = 15000
$ws.Range("O97").Value = 16000
$ws.Range("P97").Value = 15500
$ws.Range("Q97").Value = "$/caja 20 kilos"
$ws.Range("R97").Value = "Región de Coquimbo"
$ws.Range("S97").Value = 775
$ws.Range("T97").Value = 20

# Ensure the date cells keep the expected date/time number format (style carried over
# from the Insert should already match, but set it explicitly to be safe).
$ws.Range("D96").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D97").NumberFormat = "YYYY-MM-DD HH:MM:SS"
